$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the 25-Dec-2022 symbol-list refresh scraped by the GitHub Actions job.
# Column D holds prices as text (inlineStr), so we force NumberFormat to "@" (Text)
# before writing each numeric-looking string - otherwise COM auto-converts the
# assignment to a floating point number and we lose the exact textual value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '242.65'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.93'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.379'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05954'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.478'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8071'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9089'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1419'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07434'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03297'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03056'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.858'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001573'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04517'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005939'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006085'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005007'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'UpBots'
$ws.Range("C21").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.007494'
$ws.Range("E21").Value = '20UpBotsUBXTBestin24h'
$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0009797'
$ws.Range("E22").Value = '21BitKanKAN'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.00007794'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.613'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.136'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006080'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1070'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002588'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007171'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005191'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0005799'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9699'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002260'

# The NumberFormat = "@" above bumps the style index on the touched D cells.
# Restore their original (style-less) formatting by pasting the format from an
# untouched column-D cell over the whole data range - this does not touch values.
$ws.Range("D6").Copy()
$ws.Range("D2:D49").PasteSpecial(-4122)
$excel.CutCopyMode = $false
